# Update gh-pages output (杭州-漫展信息.xlsx) to the data generated at 456a3b4.
#
# Sheets: 展览 (Exhibition) / 演出 (Performance) / 本地生活 (Local Life) / 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, $row, $col, $text)
    # Force the cell to store a literal text value even when the string looks
    # like a date (e.g. "2024-11-02"), which Excel would otherwise silently
    # convert to a date serial number. Temporarily mark the cell as text,
    # assign, then restore formatting to match its (unstyled) neighbours so
    # no stray number-format sticks around on the cell.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $ws.Cells.Item($row, 3).Copy()
    $cell.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibition)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F10").Value2 = 377
$ws1.Range("F11").Value2 = 432
$ws1.Range("F13").Value2 = 319
$ws1.Range("F14").Value2 = 366
$ws1.Range("F15").Value2 = 45
$ws1.Range("F16").Value2 = 67
$ws1.Range("D17").Value2 = "杭州in77店D区B2层B2007室 三月兽mini杭州店"
$ws1.Range("F18").Value2 = 562
$ws1.Range("F19").Value2 = 1468
$ws1.Range("F20").Value2 = 5732
$ws1.Range("F21").Value2 = 93
$ws1.Range("F22").Value2 = 1613
$ws1.Range("F24").Value2 = 61
$ws1.Range("F26").Value2 = 5323
$ws1.Range("F27").Value2 = 5323
$ws1.Range("F30").Value2 = 1546
$ws1.Range("F33").Value2 = 61
$ws1.Range("F34").Value2 = 47
$ws1.Range("F36").Value2 = 108

# Insert a brand-new row 37 ("杭州·BanGDream! Only同人展"), pushing the old
# row 37 ("杭州·岚梦国潮·夏日盛会") down to row 38 and the old row 38
# ("杭州·原神X崩坏X星铁旅行盛宴·同人only首展") down to row 39.
$ws1.Rows.Item(37).Insert()

# The inserted row's "#" column (A) loses its style on insert - restore it
# from the row below (which still carries the original numbering style).
$ws1.Cells.Item(38, 1).Copy()
$ws1.Cells.Item(37, 1).PasteSpecial(-4122)

# Renumber the running index (col A, 0-based "#") for the new row and the two
# rows that shifted down, keeping the sequence contiguous (36, 37, 38 - same
# as the source data generator, which just enumerates every row).
$ws1.Range("A37").Value2 = 36
$ws1.Range("A38").Value2 = 37
$ws1.Range("A39").Value2 = 38

Set-TextValue $ws1 37 2 "2024-11-02"
$ws1.Range("C37").Value2 = "杭州·BanGDream! Only同人展"
$ws1.Range("D37").Value2 = "石祥路与丽水北路交叉口 大运河音乐公园"
$ws1.Range("E37").Value2 = "2024.11.02 10:00-11.03 20:00"
$ws1.Range("F37").Value2 = 0
$ws1.Range("G37").Value2 = 89
$ws1.Range("H37").Value2 = "https://show.bilibili.com/platform/detail.html?id=91168"
$ws1.Range("I37").Value2 = "//i0.hdslb.com/bfs/openplatform/202408/0vTxEVyz1724222524879.jpeg"

# Old row 38 (now row 39) needs its "想去人数" (F) bumped from 3810 to 3812.
$ws1.Range("F39").Value2 = 3812

# ---------------------------------------------------------------------------
# Sheet "演出" (Performance)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F4").Value2 = 39
$ws2.Range("F5").Value2 = 161
$ws2.Range("F8").Value2 = 187

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local Life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F3").Value2 = 585
$ws3.Range("F4").Value2 = 2157

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All Types)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value2 = 585
$ws4.Range("F4").Value2 = 2157
$ws4.Range("F12").Value2 = 377
$ws4.Range("F13").Value2 = 432
$ws4.Range("F14").Value2 = 319
$ws4.Range("F15").Value2 = 366
$ws4.Range("F16").Value2 = 45
$ws4.Range("F17").Value2 = 67
$ws4.Range("F21").Value2 = 1468
$ws4.Range("F22").Value2 = 5732
$ws4.Range("F23").Value2 = 93
$ws4.Range("F24").Value2 = 1613
$ws4.Range("F30").Value2 = 5323
$ws4.Range("F31").Value2 = 5323
$ws4.Range("F34").Value2 = 1546
$ws4.Range("F37").Value2 = 47
$ws4.Range("F39").Value2 = 108
$ws4.Range("F47").Value2 = 3812
